# Update countries & provincias Spain
# Refreshes the COVID country table ("Pais" sheet) with the next data pull:
#  - bump the "last updated" timestamp
#  - update totals for several countries
#  - Etiopia overtakes China, and Estado de Palestina overtakes Libia in the
#    ranking, so those two row pairs swap which country/data they hold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title: "Datos actualizados a ..." timestamp
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 21:05"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8063681
$ws.Range("C4").Value = 25892
$ws.Range("D4").Value = 5203991
$ws.Range("E4").Value = 2639258
$ws.Range("G4").Value = 421
$ws.Range("H4").Value = 220432

# Row 5: India
$ws.Range("B5").Value = 7236995
$ws.Range("C5").Value = 63430
$ws.Range("D5").Value = 6298606
$ws.Range("E5").Value = 827772
$ws.Range("G5").Value = 723
$ws.Range("H5").Value = 110617

# Row 13: Francia
$ws.Range("B13").Value = 756472
$ws.Range("C13").Value = 12993
$ws.Range("D13").Value = 102680
$ws.Range("E13").Value = 620850
$ws.Range("G13").Value = 117
$ws.Range("H13").Value = 32942

# Row 25: Alemania
$ws.Range("B25").Value = 334227
$ws.Range("C25").Value = 3133
$ws.Range("E25").Value = 45393
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 9734

# Rows 52/53: Etiopia's update pushes it ahead of China in the ranking, so
# row 52 becomes Etiopia (with fresh numbers) and row 53 becomes China
# (holding the stale numbers that used to belong to row 52).
$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 85718
$ws.Range("C52").Value = 582
$ws.Range("D52").Value = 39307
$ws.Range("E52").Value = 45106
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 1305

$ws.Range("A53").Value = "China"
$ws.Range("B53").Value = 85591
$ws.Range("C53").Value = 13
$ws.Range("D53").Value = 80729
$ws.Range("E53").Value = 228
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 4634

# Rows 70/71: Estado de Palestina's update pushes it ahead of Libia in the
# ranking, same swap pattern as above.
$ws.Range("A70").Value = "Estado de Palestina"
$ws.Range("B70").Value = 45200
$ws.Range("C70").Value = 516
$ws.Range("D70").Value = 38841
$ws.Range("E70").Value = 5968
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 391

$ws.Range("A71").Value = "Libia"
$ws.Range("B71").Value = 44985
$ws.Range("C71").Value = 1164
$ws.Range("D71").Value = 25007
$ws.Range("E71").Value = 19322
$ws.Range("G71").Value = 12
$ws.Range("H71").Value = 656

# Row 103: Namibia
$ws.Range("B103").Value = 12000
$ws.Range("C103").Value = 11
$ws.Range("D103").Value = 10120
$ws.Range("E103").Value = 1751

# Row 152: Sudan del Sur
$ws.Range("B152").Value = 2798
$ws.Range("C152").Value = 11
$ws.Range("E152").Value = 1453

# Row 164: Lesoto
$ws.Range("B164").Value = 1822
$ws.Range("C164").Value = 17
$ws.Range("E164").Value = 819

# Row 189: Monaco
$ws.Range("B189").Value = 241
$ws.Range("C189").Value = 5
$ws.Range("D189").Value = 217
$ws.Range("E189").Value = 22
